$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 11 new data rows (18-28) so the 3 worker rows (16-18) become the
#    14 period rows (16-29) for the single remaining worker.
# ---------------------------------------------------------------------------
$ws.Rows("18:28").Insert()

# Copy the formatting of row 17 (a "middle" data row) into the newly
# inserted rows 18-28 so they pick up the same borders/number formats.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Header / summary block updates
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 2565703
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 14

# ---------------------------------------------------------------------------
# 3) Fill the 14 data rows (16-29) with the remaining worker's info across
#    the 14 mora periods (2502 down to 2401).
# ---------------------------------------------------------------------------
$docType = "CC"
$docNum = "1051823100"
$name = "EVA LUZ DIAZ CATALAN"
$periods = @("2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = $docType       # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $docNum         # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $name           # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $periods[$i]    # E - Periodo Mora
    if ($row -eq 29) {
        $ws.Cells.Item($row, 6).Value = 57873
    } else {
        $ws.Cells.Item($row, 6).Value = 192910
    }
    $ws.Cells.Item($row, 7).Value = 4822731         # G - Salario Basico
}

# ---------------------------------------------------------------------------
# 4) Column D width shrinks now that the name column holds shorter text
#    (closest achievable snap to the recorded bestFit width of 21.6328125).
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 20.8
